# Update the "Data" sheet:
#  - Add the 2020 actual "Real Instantaneous Peak Load" value into F42
#  - Clear out the forecast rows 43:62 (columns A:E, G:H) that held placeholder
#    projection data, and drop the computed I:K formula/value cells entirely
#  - Move the active-cell selection to O38

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# --- F42: new actual data point (style copied from the cell above it) ---
$ws.Range("F41").Copy()
$ws.Range("F42").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("F42").Value = 46103.7

# --- Rows 43:62: clear out the stale forecast data ---
$ws.Range("A43:E62").ClearContents()
$ws.Range("G43:K62").ClearContents()

# --- Update the saved selection/active cell ---
$ws.Range("O38").Select()
